$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new row of data (row 11), matching the pattern of the
# preceding rows: Column A = title, Column B = subject, Column C = link.
$ws.Range("A11").Value = "Banco de Imagens"
$ws.Range("B11").Value = "Banco de Imagens Pedro "
$ws.Range("C11").Value = "https://www.pexels.com/"

# Update the active cell selection to C16, as recorded in the saved view.
$ws.Range("C16").Select()
